$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 8519.5
$ws.Range("I46").Value = 7138
$ws.Range("J46").Value = 8795.799999999999
$ws.Range("K46").Value = 7138
$ws.Range("L46").Value = 8795.799999999999
$ws.Range("M46").Value = -6819
$ws.Range("N46").Value = -9433.799999999999

# Row 74
$ws.Range("H74").Value = 1129.8125
$ws.Range("I74").Value = 827.4074000000001
$ws.Range("J74").Value = 2762.8
$ws.Range("K74").Value = 827.4074000000001
$ws.Range("L74").Value = 2762.8
$ws.Range("M74").Value = 46.59259999999995
$ws.Range("N74").Value = -4510.8

# Row 77
$ws.Range("H77").Value = 1129.8125
$ws.Range("I77").Value = 827.4074000000001
$ws.Range("J77").Value = 2762.8
$ws.Range("K77").Value = 4137.037
$ws.Range("L77").Value = 13814
$ws.Range("M77").Value = 230.9629999999997
$ws.Range("N77").Value = -22550

# Row 122
$ws.Range("H122").Value = 1540.5952
$ws.Range("I122").Value = 1289.2333
$ws.Range("J122").Value = 2169
$ws.Range("K122").Value = 3867.699900000001
$ws.Range("L122").Value = 6507
$ws.Range("M122").Value = -1417.699900000001
$ws.Range("N122").Value = -11407

# Row 132
$ws.Range("H132").Value = 2318.9077
$ws.Range("I132").Value = 2200.2932
$ws.Range("J132").Value = 3301.7144
$ws.Range("K132").Value = 6600.8796
$ws.Range("L132").Value = 9905.143199999999
$ws.Range("M132").Value = -4070.8796
$ws.Range("N132").Value = -14965.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 23064.25
$ws.Range("I82").Value = 5257
$ws.Range("J82").Value = 29000
$ws.Range("K82").Value = 5257
$ws.Range("L82").Value = 29000
$ws.Range("M82").Value = -4874
$ws.Range("N82").Value = -29766

# Row 85
$ws.Range("H85").Value = 23064.25
$ws.Range("I85").Value = 5257
$ws.Range("J85").Value = 29000
$ws.Range("K85").Value = 5257
$ws.Range("L85").Value = 29000
$ws.Range("M85").Value = -3931
$ws.Range("N85").Value = -31652

# Row 102
$ws.Range("H102").Value = 26100
$ws.Range("I102").Value = 15000
$ws.Range("K102").Value = 15000
$ws.Range("M102").Value = -11755

# Row 134
$ws.Range("H134").Value = 2859.725
$ws.Range("I134").Value = 2929.3242
$ws.Range("J134").Value = 2001.3334
$ws.Range("K134").Value = 8787.972600000001
$ws.Range("L134").Value = 6004.0002
$ws.Range("M134").Value = -6252.972600000001
$ws.Range("N134").Value = -11074.0002

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2157.0984
$ws.Range("I31").Value = 1314.7778
$ws.Range("J31").Value = 4526.125
$ws.Range("K31").Value = 1314.7778
$ws.Range("L31").Value = 4526.125
$ws.Range("M31").Value = -1019.7778
$ws.Range("N31").Value = -5116.125

# Row 34
$ws.Range("H34").Value = 2157.0984
$ws.Range("I34").Value = 1314.7778
$ws.Range("J34").Value = 4526.125
$ws.Range("K34").Value = 1314.7778
$ws.Range("L34").Value = 4526.125
$ws.Range("M34").Value = -1112.7778
$ws.Range("N34").Value = -4930.125

# Row 58
$ws.Range("H58").Value = 13516930
$ws.Range("I58").Value = 2684.1667
$ws.Range("J58").Value = 26319900
$ws.Range("K58").Value = 2684.1667
$ws.Range("L58").Value = 26319900
$ws.Range("M58").Value = -2481.1667
$ws.Range("N58").Value = -26320306

# Row 105
$ws.Range("H105").Value = 1912.2222
$ws.Range("I105").Value = 2172.8572
$ws.Range("K105").Value = 2172.8572
$ws.Range("M105").Value = -425.8571999999999

# Row 132
$ws.Range("H132").Value = 1885.5
$ws.Range("I132").Value = 1767.65
$ws.Range("J132").Value = 2474.75
$ws.Range("K132").Value = 5302.950000000001
$ws.Range("L132").Value = 7424.25
$ws.Range("M132").Value = -2772.950000000001
$ws.Range("N132").Value = -12484.25

# Row 134
$ws.Range("H134").Value = 5258.3184
$ws.Range("I134").Value = 5371.278
$ws.Range("J134").Value = 4750
$ws.Range("K134").Value = 16113.834
$ws.Range("L134").Value = 14250
$ws.Range("M134").Value = -13578.834
$ws.Range("N134").Value = -19320

# Row 136
$ws.Range("H136").Value = 13516930
$ws.Range("I136").Value = 2684.1667
$ws.Range("J136").Value = 26319900
$ws.Range("K136").Value = 8052.500100000001
$ws.Range("L136").Value = 78959700
$ws.Range("M136").Value = -5502.500100000001
$ws.Range("N136").Value = -78964800

$ws = $wb.Worksheets.Item("CUL")
# Row 42
$ws.Range("H42").Value = 4202
$ws.Range("J42").Value = 4202
$ws.Range("L42").Value = 12606
$ws.Range("N42").Value = -13674

# Row 131
$ws.Range("H131").Value = 1561.7567
$ws.Range("I131").Value = 3033.625
$ws.Range("J131").Value = 1155.7241
$ws.Range("K131").Value = 9100.875
$ws.Range("L131").Value = 3467.1723
$ws.Range("M131").Value = -4060.875
$ws.Range("N131").Value = -13547.1723

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 211433.58
$ws.Range("I18").Value = 535002.5
$ws.Range("J18").Value = 82006
$ws.Range("K18").Value = 535002.5
$ws.Range("L18").Value = 82006
$ws.Range("M18").Value = -534709.5
$ws.Range("N18").Value = -82592

# Row 122
$ws.Range("H122").Value = 5060.722
$ws.Range("I122").Value = 4706.4287
$ws.Range("J122").Value = 5286.1816
$ws.Range("K122").Value = 14119.2861
$ws.Range("L122").Value = 15858.5448
$ws.Range("M122").Value = -11669.2861
$ws.Range("N122").Value = -20758.5448

# Row 132
$ws.Range("H132").Value = 2745.3389
$ws.Range("I132").Value = 2552.0527
$ws.Range("J132").Value = 3095.0952
$ws.Range("K132").Value = 7656.158100000001
$ws.Range("L132").Value = 9285.285600000001
$ws.Range("M132").Value = -5126.158100000001
$ws.Range("N132").Value = -14345.2856

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1657.75
$ws.Range("I40").Value = 1588.1111
$ws.Range("J40").Value = 1866.6666
$ws.Range("K40").Value = 1588.1111
$ws.Range("L40").Value = 1866.6666
$ws.Range("M40").Value = -1452.1111
$ws.Range("N40").Value = -2138.6666

# Row 61
$ws.Range("H61").Value = 1924.1428
$ws.Range("I61").Value = 578.1667
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 578.1667
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -376.1667
$ws.Range("N61").Value = -10404

# Row 113
$ws.Range("H113").Value = 1924.1428
$ws.Range("I113").Value = 578.1667
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 578.1667
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 1591.8333
$ws.Range("N113").Value = -14340

# Row 122
$ws.Range("H122").Value = 2562.5557
$ws.Range("I122").Value = 2548.1538
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 7644.4614
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -5194.4614
$ws.Range("N122").Value = -12700

# Row 132
$ws.Range("H132").Value = 5940.1875
$ws.Range("I132").Value = 2411.2727
$ws.Range("K132").Value = 7233.8181
$ws.Range("M132").Value = -4703.8181

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 80005
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 122
$ws.Range("H122").Value = 3206.1765
$ws.Range("I122").Value = 2600
$ws.Range("J122").Value = 3745
$ws.Range("K122").Value = 7800
$ws.Range("L122").Value = 11235
$ws.Range("M122").Value = -5350
$ws.Range("N122").Value = -16135
